$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: C2 model text changes from "ONIX 1.2 L/20" to "ONIX 1.2 L/25"
$ws.Range("C2").Value = "ONIX 1.2 L/25"

# Row 2: replace external-link formulas (F2/G2/H2) with plain literal values
# (order matters so new shared strings land in the same order as the target)
$ws.Range("G2").Value = "ZZZ111BB00"
$ws.Range("H2").Value = "ZZZ111BB0011"
$ws.Range("F2").Value = "ZZZ113"

# Row 3: model text now also "ONIX 1.2 L/25", SumaAsegurada bumped by 1
$ws.Range("C3").Value = "ONIX 1.2 L/25"
$ws.Range("D3").Value = 1700001

# Row 3: replace external-link formulas (F3/G3/H3) with plain literal values
$ws.Range("F3").Value = "ZZZ112"
$ws.Range("G3").Value = "ZZZ111BB00"
$ws.Range("H3").Value = "ZZZ111BB0011"

# Remove the now-unused external reference to the source workbook entirely
$sources = $wb.LinkSources()
if ($sources) {
    foreach ($s in $sources) {
        $wb.BreakLink($s, 1)
    }
}

# Move the active selection from F3 to H6
$ws.Range("H6").Select()
